$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.412.09'
$ws.Range("E2").Value = '  +0.30%  '
$ws.Range("D3").Value = '1.869.52'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.14'
$ws.Range("E5").Value = '  +0.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7046'
$ws.Range("E6").Value = '  -1.02%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07944'
$ws.Range("E8").Value = '  -1.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3141'
$ws.Range("E9").Value = '  -0.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.56'
$ws.Range("E10").Value = '  -1.76%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07859'
$ws.Range("E11").Value = '  -4.66%  '
$ws.Range("D12").Value = '1.868.13'
$ws.Range("E12").Value = '  -1.77%  '
$ws.Range("B13").Value = 'Litecoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '93.87'
$ws.Range("E13").Value = '  -0.94%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.192'
$ws.Range("E14").Value = '  -1.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7032'
$ws.Range("E15").Value = '  -1.33%  '
$ws.Range("E16").Value = '  +2.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008416'
$ws.Range("E17").Value = '  -1.65%  '
$ws.Range("D18").Value = '29.405.39'
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '253.88'
$ws.Range("E19").Value = '  +4.11%  '
$ws.Range("D20").Value = '2.129.39'
$ws.Range("E20").Value = '  -0.96%  '
$ws.Range("E21").Value = '  -1.03%  '
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.648'
$ws.Range("E23").Value = '  -1.74%  '
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.017'
$ws.Range("E26").Value = '  -0.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.19'
$ws.Range("E27").Value = '  -0.81%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.82'
$ws.Range("E28").Value = '  +1.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.505'
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.330'
$ws.Range("E30").Value = '  -2.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.265'
$ws.Range("E31").Value = '  -1.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.214'
$ws.Range("E32").Value = '  +2.58%  '
$ws.Range("E33").Value = '  -1.49%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.896'
$ws.Range("E34").Value = '  -2.35%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7513'
$ws.Range("E35").Value = '  -1.81%  '
$ws.Range("E36").Value = '  -0.37%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.711'
$ws.Range("E37").Value = '  +0.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01888'
$ws.Range("E38").Value = '  +0.55%  '
$ws.Range("D39").Value = '1.288.16'
$ws.Range("E39").Value = '  +2.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.766'
$ws.Range("E40").Value = '  +0.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8948'
$ws.Range("E41").Value = '  -2.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '109.22'
$ws.Range("E42").Value = '  -3.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.029'
$ws.Range("E43").Value = '  -7.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '71.23'
$ws.Range("E44").Value = '  -3.99%  '
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000127'
$ws.Range("E46").Value = '  -5.06%  '
$ws.Range("D47").Value = '2.027.75'
$ws.Range("E47").Value = '  -0.92%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.800'
$ws.Range("E48").Value = '  -0.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.586'
$ws.Range("E49").Value = '  +1.28%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.5178'
$ws.Range("E50").Value = '  -0.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4313'
$ws.Range("E51").Value = '  -0.98%  '
